# NOTE: row numbers used below are the numbers *before* the later
# row-10 insertion (which will shift everything from row 10 downward
# by one row), chosen so that after the insert they land on the
# final target row numbers (16 and 17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: add the new "HAR+ OSAKA" row (will become row 17 after the insert below) ---
$ws.Range("A16").Value = "HAR+ OSAKA"
$ws.Range("B16").Value = "prediction accuracy on test set: 83.3333%"
$ws.Range("C16").Value = 1927064
$ws.Range("D16").Value = 7708
$ws.Range("E16").Value = "Wavelet + RF"
$ws.Rows("16:16").RowHeight = 54

# Blank separator row just above it (will become row 16 after the insert below)
$ws.Rows("15:15").RowHeight = 34.2

# --- Step 2: rename the "Three Sensors at Waist" rows (still rows 13 & 14 at this point)
#     to "One sensor at Waist" ---
$ws.Range("A13").Value = "One sensor at Waist"
$ws.Range("A14").Value = "One sensor at Waist"

# --- Step 3: insert a new row at position 10, shifting the renamed rows down to 14 & 15 ---
$ws.Rows("10:10").Insert()

# --- Step 4: fill the newly inserted row 10 with the reduced-sample "Waist" results ---
$ws.Range("A10").Value = "Waist"
$ws.Range("B10").Value = "Reduced number of samples to be same as osaka dataset : 95.089%"
$ws.Range("C10").Value = 702416
$ws.Range("D10").Value = 1404
$ws.Range("E10").Value = "Wavelet + RF"

# Clear whatever formatting Insert() copied down, then apply the correct styling
$ws.Range("A10:E10").ClearFormats()
$ws.Range("B10").WrapText = $true
$ws.Range("B10").NumberFormat = "0.00%"
$ws.Rows("10:10").RowHeight = 66.6

# --- Step 5: update the view selection ---
$ws.Range("F10").Select()
